# Updates cryptos list prices / volume percentages (and the swapped
# Polkadot / WrappedBTC rows) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.212.65'
$ws.Range('E2').Value = '  +2.81%  '
$ws.Range('D3').Value = '3.043.63'
$ws.Range('E3').Value = '  +1.70%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'596.22"
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').Value = "'154.97"
$ws.Range('E6').Value = '  +7.57%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.040.36'
$ws.Range('E8').Value = '  +1.58%  '
$ws.Range('D9').Value = "'0.517"
$ws.Range('E9').Value = '  -0.34%  '
$ws.Range('D10').Value = "'6.88"
$ws.Range('E10').Value = '  +13.59%  '
$ws.Range('E11').Value = '  +3.64%  '
$ws.Range('D12').Value = "'0.466"
$ws.Range('E12').Value = '  +2.32%  '
$ws.Range('E13').Value = '  +2.98%  '
$ws.Range('D14').Value = "'35.90"
$ws.Range('E14').Value = '  +4.54%  '
$ws.Range('E15').Value = '  +2.02%  '
$ws.Range('D16').Value = '3.548.95'
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '63.210.11'
$ws.Range('E17').Value = '  +2.75%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = "'7.09"
$ws.Range('E18').Value = '  +1.89%  '
$ws.Range('D19').Value = '3.046.95'
$ws.Range('E19').Value = '  +1.79%  '
$ws.Range('D20').Value = "'455.70"
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').Value = "'14.34"
$ws.Range('E21').Value = '  +2.49%  '
$ws.Range('E22').Value = '  +2.35%  '
$ws.Range('D23').Value = "'7.55"
$ws.Range('E23').Value = '  +3.12%  '
$ws.Range('D24').Value = "'83.15"
$ws.Range('E24').Value = '  +2.00%  '
$ws.Range('D25').Value = "'11.26"
$ws.Range('E25').Value = '  +4.20%  '
$ws.Range('D26').Value = "'2.31"
$ws.Range('E26').Value = '  +3.93%  '
$ws.Range('E27').Value = '  +3.80%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = "'7.51"
$ws.Range('E29').Value = '  +4.40%  '
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('E31').Value = '  +9.11%  '
$ws.Range('D32').Value = "'1.00"
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('E33').Value = '  +1.41%  '
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('D35').Value = '0.0₃0865'
$ws.Range('E35').Value = '  +4.32%  '
$ws.Range('D36').Value = "'1.04"
$ws.Range('D37').Value = "'5.96"
$ws.Range('E37').Value = '  +3.14%  '
$ws.Range('D38').Value = "'3.22"
$ws.Range('E38').Value = '  +12.32%  '
$ws.Range('E39').Value = '  +3.09%  '
$ws.Range('E40').Value = '  +5.84%  '
$ws.Range('D41').Value = "'50.44"
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('E42').Value = '  -0.70%  '
$ws.Range('E43').Value = '  +12.87%  '
$ws.Range('D44').Value = "'43.55"
$ws.Range('E44').Value = '  +9.09%  '
$ws.Range('D45').Value = "'397.61"
$ws.Range('E46').Value = '  +2.80%  '
$ws.Range('D47').Value = '2.730.58'
$ws.Range('E47').Value = '  +1.52%  '
$ws.Range('D48').Value = "'132.40"
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('E49').Value = '  +7.40%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').Value = "'24.54"
$ws.Range('E51').Value = '  +4.40%  '
